$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hgf"
$ws.Cells.Item(2, 3).Value = "Met"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1088396666666667
$ws.Cells.Item(2, 8).Value = 0.326519
$ws.Cells.Item(2, 9).Value = 0.002750770615347974
$ws.Cells.Item(2, 10).Value = 0.002750770615347974
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.167735333333333
$ws.Cells.Item(2, 14).Value = 9.503206
$ws.Cells.Item(2, 15).Value = 0.1182666224938439
$ws.Cells.Item(2, 16).Value = 0.1182666224938439
$ws.Cells.Item(2, 17).Value = 0.3447752577682222
$ws.Cells.Item(2, 18).Value = 3.102977319914
$ws.Cells.Item(2, 19).Value = 0.0003253243499325175
$ws.Cells.Item(2, 20).Value = 0.0003253243499325176

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hgf"
$ws.Cells.Item(3, 3).Value = "Met"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1088396666666667
$ws.Cells.Item(3, 8).Value = 0.326519
$ws.Cells.Item(3, 9).Value = 0.002750770615347974
$ws.Cells.Item(3, 10).Value = 0.002750770615347974
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.9421210000000001
$ws.Cells.Item(3, 14).Value = 2.826363
$ws.Cells.Item(3, 15).Value = 0.03517385669126484
$ws.Cells.Item(3, 16).Value = 0.03517385669126484
$ws.Cells.Item(3, 17).Value = 0.1025401355996667
$ws.Cells.Item(3, 18).Value = 0.922861220397
$ws.Cells.Item(3, 19).Value = [double]"9.675521141479202E-05"
$ws.Cells.Item(3, 20).Value = [double]"9.675521141479204E-05"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hgf"
$ws.Cells.Item(4, 3).Value = "Met"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1088396666666667
$ws.Cells.Item(4, 8).Value = 0.326519
$ws.Cells.Item(4, 9).Value = 0.002750770615347974
$ws.Cells.Item(4, 10).Value = 0.002750770615347974
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.442127
$ws.Cells.Item(4, 14).Value = 4.326381
$ws.Cells.Item(4, 15).Value = 0.05384145818700961
$ws.Cells.Item(4, 16).Value = 0.0538414581870096
$ws.Cells.Item(4, 17).Value = 0.156960621971
$ws.Cells.Item(4, 18).Value = 1.412645597739
$ws.Cells.Item(4, 19).Value = 0.0001481055010683126
$ws.Cells.Item(4, 20).Value = 0.0001481055010683126

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Hgf"
$ws.Cells.Item(5, 3).Value = "Met"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.1088396666666667
$ws.Cells.Item(5, 8).Value = 0.326519
$ws.Cells.Item(5, 9).Value = 0.002750770615347974
$ws.Cells.Item(5, 10).Value = 0.002750770615347974
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 21.232711
$ws.Cells.Item(5, 14).Value = 63.69813300000001
$ws.Cells.Item(5, 15).Value = 0.7927180626278817
$ws.Cells.Item(5, 16).Value = 0.7927180626278817
$ws.Cells.Item(5, 17).Value = 2.310961187669667
$ws.Cells.Item(5, 18).Value = 20.798650689027
$ws.Cells.Item(5, 19).Value = 0.002180585552932351
$ws.Cells.Item(5, 20).Value = 0.002180585552932352

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hgf"
$ws.Cells.Item(6, 3).Value = "Met"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 15.07419333333333
$ws.Cells.Item(6, 8).Value = 45.22258
$ws.Cells.Item(6, 9).Value = 0.3809791902285103
$ws.Cells.Item(6, 10).Value = 0.3809791902285103
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.167735333333333
$ws.Cells.Item(6, 14).Value = 9.503206
$ws.Cells.Item(6, 15).Value = 0.1182666224938439
$ws.Cells.Item(6, 16).Value = 0.1182666224938439
$ws.Cells.Item(6, 17).Value = 47.75105484349778
$ws.Cells.Item(6, 18).Value = 429.75949359148
$ws.Cells.Item(6, 19).Value = 0.04505712206876557
$ws.Cells.Item(6, 20).Value = 0.04505712206876558

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hgf"
$ws.Cells.Item(7, 3).Value = "Met"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 15.07419333333333
$ws.Cells.Item(7, 8).Value = 45.22258
$ws.Cells.Item(7, 9).Value = 0.3809791902285103
$ws.Cells.Item(7, 10).Value = 0.3809791902285103
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9421210000000001
$ws.Cells.Item(7, 14).Value = 2.826363
$ws.Cells.Item(7, 15).Value = 0.03517385669126484
$ws.Cells.Item(7, 16).Value = 0.03517385669126484
$ws.Cells.Item(7, 17).Value = 14.20171409739334
$ws.Cells.Item(7, 18).Value = 127.81542687654
$ws.Cells.Item(7, 19).Value = 0.01340050743945175
$ws.Cells.Item(7, 20).Value = 0.01340050743945175

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Hgf"
$ws.Cells.Item(8, 3).Value = "Met"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.07419333333333
$ws.Cells.Item(8, 8).Value = 45.22258
$ws.Cells.Item(8, 9).Value = 0.3809791902285103
$ws.Cells.Item(8, 10).Value = 0.3809791902285103
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.442127
$ws.Cells.Item(8, 14).Value = 4.326381
$ws.Cells.Item(8, 15).Value = 0.05384145818700961
$ws.Cells.Item(8, 16).Value = 0.0538414581870096
$ws.Cells.Item(8, 17).Value = 21.73890120922
$ws.Cells.Item(8, 18).Value = 195.65011088298
$ws.Cells.Item(8, 19).Value = 0.02051247514080912
$ws.Cells.Item(8, 20).Value = 0.02051247514080912

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Hgf"
$ws.Cells.Item(9, 3).Value = "Met"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.07419333333333
$ws.Cells.Item(9, 8).Value = 45.22258
$ws.Cells.Item(9, 9).Value = 0.3809791902285103
$ws.Cells.Item(9, 10).Value = 0.3809791902285103
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 21.232711
$ws.Cells.Item(9, 14).Value = 63.69813300000001
$ws.Cells.Item(9, 15).Value = 0.7927180626278817
$ws.Cells.Item(9, 16).Value = 0.7927180626278817
$ws.Cells.Item(9, 17).Value = 320.0659906047933
$ws.Cells.Item(9, 18).Value = 2880.59391544314
$ws.Cells.Item(9, 19).Value = 0.3020090855794839
$ws.Cells.Item(9, 20).Value = 0.3020090855794839

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Hgf"
$ws.Cells.Item(10, 3).Value = "Met"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 23.69325166666667
$ws.Cells.Item(10, 8).Value = 71.07975500000001
$ws.Cells.Item(10, 9).Value = 0.5988138558556568
$ws.Cells.Item(10, 10).Value = 0.5988138558556569
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.167735333333333
$ws.Cells.Item(10, 14).Value = 9.503206
$ws.Cells.Item(10, 15).Value = 0.1182666224938439
$ws.Cells.Item(10, 16).Value = 0.1182666224938439
$ws.Cells.Item(10, 17).Value = 75.05395046605889
$ws.Cells.Item(10, 18).Value = 675.4855541945301
$ws.Cells.Item(10, 19).Value = 0.07081969223456402
$ws.Cells.Item(10, 20).Value = 0.07081969223456404

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Hgf"
$ws.Cells.Item(11, 3).Value = "Met"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 23.69325166666667
$ws.Cells.Item(11, 8).Value = 71.07975500000001
$ws.Cells.Item(11, 9).Value = 0.5988138558556568
$ws.Cells.Item(11, 10).Value = 0.5988138558556569
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.9421210000000001
$ws.Cells.Item(11, 14).Value = 2.826363
$ws.Cells.Item(11, 15).Value = 0.03517385669126484
$ws.Cells.Item(11, 16).Value = 0.03517385669126484
$ws.Cells.Item(11, 17).Value = 22.32190995345167
$ws.Cells.Item(11, 18).Value = 200.897189581065
$ws.Cells.Item(11, 19).Value = 0.0210625927506106
$ws.Cells.Item(11, 20).Value = 0.0210625927506106

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Hgf"
$ws.Cells.Item(12, 3).Value = "Met"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 23.69325166666667
$ws.Cells.Item(12, 8).Value = 71.07975500000001
$ws.Cells.Item(12, 9).Value = 0.5988138558556568
$ws.Cells.Item(12, 10).Value = 0.5988138558556569
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.442127
$ws.Cells.Item(12, 14).Value = 4.326381
$ws.Cells.Item(12, 15).Value = 0.05384145818700961
$ws.Cells.Item(12, 16).Value = 0.0538414581870096
$ws.Cells.Item(12, 17).Value = 34.168677946295
$ws.Cells.Item(12, 18).Value = 307.518101516655
$ws.Cells.Item(12, 19).Value = 0.03224101118185434
$ws.Cells.Item(12, 20).Value = 0.03224101118185435

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Hgf"
$ws.Cells.Item(13, 3).Value = "Met"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 23.69325166666667
$ws.Cells.Item(13, 8).Value = 71.07975500000001
$ws.Cells.Item(13, 9).Value = 0.5988138558556568
$ws.Cells.Item(13, 10).Value = 0.5988138558556569
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 21.232711
$ws.Cells.Item(13, 14).Value = 63.69813300000001
$ws.Cells.Item(13, 15).Value = 0.7927180626278817
$ws.Cells.Item(13, 16).Value = 0.7927180626278817
$ws.Cells.Item(13, 17).Value = 503.0719652886017
$ws.Cells.Item(13, 18).Value = 4527.647687597416
$ws.Cells.Item(13, 19).Value = 0.4746905596886279
$ws.Cells.Item(13, 20).Value = 0.474690559688628

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Hgf"
$ws.Cells.Item(14, 3).Value = "Met"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.6906883333333335
$ws.Cells.Item(14, 8).Value = 2.072065
$ws.Cells.Item(14, 9).Value = 0.01745618330048481
$ws.Cells.Item(14, 10).Value = 0.01745618330048481
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.167735333333333
$ws.Cells.Item(14, 14).Value = 9.503206
$ws.Cells.Item(14, 15).Value = 0.1182666224938439
$ws.Cells.Item(14, 16).Value = 0.1182666224938439
$ws.Cells.Item(14, 17).Value = 2.187917837821112
$ws.Cells.Item(14, 18).Value = 19.69126054039
$ws.Cells.Item(14, 19).Value = 0.002064483840581779
$ws.Cells.Item(14, 20).Value = 0.002064483840581779

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Hgf"
$ws.Cells.Item(15, 3).Value = "Met"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.6906883333333335
$ws.Cells.Item(15, 8).Value = 2.072065
$ws.Cells.Item(15, 9).Value = 0.01745618330048481
$ws.Cells.Item(15, 10).Value = 0.01745618330048481
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.9421210000000001
$ws.Cells.Item(15, 14).Value = 2.826363
$ws.Cells.Item(15, 15).Value = 0.03517385669126484
$ws.Cells.Item(15, 16).Value = 0.03517385669126484
$ws.Cells.Item(15, 17).Value = 0.6507119832883336
$ws.Cells.Item(15, 18).Value = 5.856407849595001
$ws.Cells.Item(15, 19).Value = 0.0006140012897877033
$ws.Cells.Item(15, 20).Value = 0.0006140012897877033

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Hgf"
$ws.Cells.Item(16, 3).Value = "Met"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.6906883333333335
$ws.Cells.Item(16, 8).Value = 2.072065
$ws.Cells.Item(16, 9).Value = 0.01745618330048481
$ws.Cells.Item(16, 10).Value = 0.01745618330048481
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.442127
$ws.Cells.Item(16, 14).Value = 4.326381
$ws.Cells.Item(16, 15).Value = 0.05384145818700961
$ws.Cells.Item(16, 16).Value = 0.0538414581870096
$ws.Cells.Item(16, 17).Value = 0.9960602940850002
$ws.Cells.Item(16, 18).Value = 8.964542646765
$ws.Cells.Item(16, 19).Value = 0.0009398663632778283
$ws.Cells.Item(16, 20).Value = 0.0009398663632778282

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Hgf"
$ws.Cells.Item(17, 3).Value = "Met"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.6906883333333335
$ws.Cells.Item(17, 8).Value = 2.072065
$ws.Cells.Item(17, 9).Value = 0.01745618330048481
$ws.Cells.Item(17, 10).Value = 0.01745618330048481
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 21.232711
$ws.Cells.Item(17, 14).Value = 63.69813300000001
$ws.Cells.Item(17, 15).Value = 0.7927180626278817
$ws.Cells.Item(17, 16).Value = 0.7927180626278817
$ws.Cells.Item(17, 17).Value = 14.66518577273834
$ws.Cells.Item(17, 18).Value = 131.986671954645
$ws.Cells.Item(17, 19).Value = 0.0138378318068375
$ws.Cells.Item(17, 20).Value = 0.0138378318068375
